# Commit: "incorporacion de srcipt para incorporar ids automaticos a los
# productos en el alimentos.html + correcciones varias"
#
# The underlying data edit removes the obsolete "Purina-Gato-Pro-Plan"
# product row (brand "Purina Pro Plan", price 0.0) from row 59 of the
# price list. Deleting the whole row shifts every following row up by
# one, which matches the diff (mergeCells ranges and shared-string
# indices all shift accordingly once the row, and its now-unused
# strings, are gone).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 59 = "Purina-Gato-Pro-Plan" / "Purina Pro Plan" / 0.0
$ws.Rows.Item(59).Delete()

$wb.Save()
